$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $text) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.ClearFormats()
}

# Row 7
Set-TextCell "B7" "20201015"
Set-TextCell "C7" "20201008-ZSAC-0001"
Set-TextCell "D7" "WILLIS TOWERS WATSON INSURANCE SERVICES WEST, "
Set-TextCell "F7" "20201009"
$ws.Range("G7").Value = 216.83

# Row 8
Set-TextCell "B8" "20201015"
Set-TextCell "C8" "20201008-ZSAC-0002"
Set-TextCell "D8" "CELEBRATION TOURS & TRAVEL"
Set-TextCell "F8" "20201013"
$ws.Range("G8").Value = 14650

# Row 9
Set-TextCell "B9" "20201015"
Set-TextCell "C9" "20201008-ZSAC-0003"
Set-TextCell "D9" "MULTI SERVICE AVIATION"
Set-TextCell "F9" "20201014"
$ws.Range("G9").Value = 6185.47

# Match the final selection recorded in the saved view state
$ws.Range("N26").Select()
